$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 59445
$ws.Range("I33").Value = 111535.664
$ws.Range("J33").Value = 843
$ws.Range("K33").Value = 111535.664
$ws.Range("L33").Value = 843
$ws.Range("M33").Value = -111306.664
$ws.Range("N33").Value = -1301
$ws.Range("H58").Value = 2834.6365
$ws.Range("I58").Value = 1432.8
$ws.Range("J58").Value = 4002.8333
$ws.Range("K58").Value = 4298.4
$ws.Range("L58").Value = 12008.4999
$ws.Range("M58").Value = -4148.4
$ws.Range("N58").Value = -12308.4999
$ws.Range("H70").Value = 1656.4286
$ws.Range("I70").Value = 1738.4
$ws.Range("J70").Value = 1451.5
$ws.Range("K70").Value = 5215.200000000001
$ws.Range("L70").Value = 4354.5
$ws.Range("M70").Value = -4945.200000000001
$ws.Range("N70").Value = -4894.5
$ws.Range("H73").Value = 1656.4286
$ws.Range("I73").Value = 1738.4
$ws.Range("J73").Value = 1451.5
$ws.Range("K73").Value = 5215.200000000001
$ws.Range("L73").Value = 4354.5
$ws.Range("M73").Value = -4279.200000000001
$ws.Range("N73").Value = -6226.5
$ws.Range("H92").Value = 216.81818
$ws.Range("I92").Value = 216.81818
$ws.Range("K92").Value = 216.81818
$ws.Range("M92").Value = 1031.18182
$ws.Range("H125").Value = 868.1539
$ws.Range("I125").Value = 800
$ws.Range("K125").Value = 7200
$ws.Range("M125").Value = -4740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 169.78572
$ws.Range("I5").Value = 182.25
$ws.Range("K5").Value = 182.25
$ws.Range("M5").Value = -70.25
$ws.Range("H45").Value = 1367.75
$ws.Range("I45").Value = 2499
$ws.Range("K45").Value = 2499
$ws.Range("M45").Value = -2122
$ws.Range("H63").Value = 11442.667
$ws.Range("J63").Value = 13855.571
$ws.Range("L63").Value = 13855.571
$ws.Range("N63").Value = -15227.571
$ws.Range("H66").Value = 11442.667
$ws.Range("J66").Value = 13855.571
$ws.Range("L66").Value = 69277.855
$ws.Range("N66").Value = -76141.855
$ws.Range("H74").Value = 5147431
$ws.Range("I74").Value = 2649657.8
$ws.Range("K74").Value = 2649657.8
$ws.Range("M74").Value = -2648783.8
$ws.Range("H77").Value = 5147431
$ws.Range("I77").Value = 2649657.8
$ws.Range("K77").Value = 13248289
$ws.Range("M77").Value = -13243921
$ws.Range("H110").Value = 8411.5
$ws.Range("I110").Value = 5700
$ws.Range("K110").Value = 5700
$ws.Range("M110").Value = -3655
$ws.Range("H132").Value = 20001216
$ws.Range("I132").Value = 1182.2
$ws.Range("J132").Value = 100001350
$ws.Range("K132").Value = 3546.6
$ws.Range("L132").Value = 300004050
$ws.Range("M132").Value = -1016.6
$ws.Range("N132").Value = -300009110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 169.78572
$ws.Range("I4").Value = 182.25
$ws.Range("K4").Value = 182.25
$ws.Range("M4").Value = -67.25
$ws.Range("H22").Value = 364.2857
$ws.Range("I22").Value = 341.66666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 341.66666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -168.66666
$ws.Range("N22").Value = -846
$ws.Range("H97").Value = 4499
$ws.Range("I97").Value = 4499
$ws.Range("K97").Value = 4499
$ws.Range("M97").Value = -3508

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1573
$ws.Range("I22").Value = 383.22223
$ws.Range("J22").Value = 4250
$ws.Range("K22").Value = 383.22223
$ws.Range("L22").Value = 4250
$ws.Range("M22").Value = -33.22223000000002
$ws.Range("N22").Value = -4950
$ws.Range("H58").Value = 1818.4166
$ws.Range("I58").Value = 1795.6875
$ws.Range("J58").Value = 1863.875
$ws.Range("K58").Value = 1795.6875
$ws.Range("L58").Value = 1863.875
$ws.Range("M58").Value = -1592.6875
$ws.Range("N58").Value = -2269.875
$ws.Range("H94").Value = 983.7143
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 983.7143
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 983.7143
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1885.7143
$ws.Range("H136").Value = 1818.4166
$ws.Range("I136").Value = 1795.6875
$ws.Range("J136").Value = 1863.875
$ws.Range("K136").Value = 5387.0625
$ws.Range("L136").Value = 5591.625
$ws.Range("M136").Value = -2837.0625
$ws.Range("N136").Value = -10691.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 203.36842
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 203.36842
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 1220.21052
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -1786.21052
$ws.Range("H51").Value = 1480.2
$ws.Range("I51").Value = 1476.5
$ws.Range("K51").Value = 4429.5
$ws.Range("M51").Value = -3969.5
$ws.Range("H54").Value = 13084.25
$ws.Range("I54").Value = 900
$ws.Range("J54").Value = 14824.857
$ws.Range("K54").Value = 2700
$ws.Range("L54").Value = 44474.571
$ws.Range("M54").Value = -2141
$ws.Range("N54").Value = -45592.571
$ws.Range("H107").Value = 566.875
$ws.Range("J107").Value = 556
$ws.Range("L107").Value = 1668
$ws.Range("N107").Value = -5508
$ws.Range("H122").Value = 1544
$ws.Range("J122").Value = 2117.1667
$ws.Range("L122").Value = 19054.5003
$ws.Range("N122").Value = -23954.5003
$ws.Range("H139").Value = 2720.182
$ws.Range("I139").Value = 2553.75
$ws.Range("K139").Value = 7661.25
$ws.Range("M139").Value = -2521.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3838.1667
$ws.Range("I122").Value = 4356.25
$ws.Range("K122").Value = 13068.75
$ws.Range("M122").Value = -10618.75
$ws.Range("H134").Value = 262496.5
$ws.Range("J134").Value = 262496.5
$ws.Range("L134").Value = 787489.5
$ws.Range("N134").Value = -792559.5
$ws.Range("H135").Value = 295158
$ws.Range("J135").Value = 295158
$ws.Range("L135").Value = 295158
$ws.Range("N135").Value = -305298

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 30900
$ws.Range("I4").Value = 30900
$ws.Range("K4").Value = 30900
$ws.Range("M4").Value = -30787
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H25").Value = 16286.667
$ws.Range("J25").Value = 17000
$ws.Range("L25").Value = 17000
$ws.Range("N25").Value = -17460
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H28").Value = 30900
$ws.Range("I28").Value = 30900
$ws.Range("K28").Value = 30900
$ws.Range("M28").Value = -30668
$ws.Range("H37").Value = 30900
$ws.Range("I37").Value = 30900
$ws.Range("K37").Value = 30900
$ws.Range("M37").Value = -30793
$ws.Range("H61").Value = 3053.0557
$ws.Range("I61").Value = 2434.8125
$ws.Range("K61").Value = 2434.8125
$ws.Range("M61").Value = -2232.8125
$ws.Range("H93").Value = 1992.85
$ws.Range("I93").Value = 1553.375
$ws.Range("K93").Value = 1553.375
$ws.Range("M93").Value = -305.375
$ws.Range("H113").Value = 3053.0557
$ws.Range("I113").Value = 2434.8125
$ws.Range("K113").Value = 2434.8125
$ws.Range("M113").Value = -264.8125
$ws.Range("H132").Value = 3306.3
$ws.Range("I132").Value = 3062.5557
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 9187.667099999999
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -6657.667099999999
$ws.Range("N132").Value = -21560
$ws.Range("H136").Value = 41669170
$ws.Range("I136").Value = 2296.2778
$ws.Range("J136").Value = 166669780
$ws.Range("K136").Value = 6888.8334
$ws.Range("L136").Value = 500009340
$ws.Range("M136").Value = -4338.8334
$ws.Range("N136").Value = -500014440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 87850
$ws.Range("J41").Value = 101719.125
$ws.Range("L41").Value = 101719.125
$ws.Range("N41").Value = -102499.125
$ws.Range("H62").Value = 3287
$ws.Range("I62").Value = 2222
$ws.Range("K62").Value = 2222
$ws.Range("M62").Value = -1598
$ws.Range("H65").Value = 3287
$ws.Range("I65").Value = 2222
$ws.Range("K65").Value = 11110
$ws.Range("M65").Value = -7990
$ws.Range("H70").Value = 30282.75
$ws.Range("J70").Value = 30282.75
$ws.Range("L70").Value = 30282.75
$ws.Range("N70").Value = -30912.75
$ws.Range("H73").Value = 30282.75
$ws.Range("J73").Value = 30282.75
$ws.Range("L73").Value = 30282.75
$ws.Range("N73").Value = -32466.75
$ws.Range("H96").Value = 5357.4287
$ws.Range("I96").Value = 7875.75
$ws.Range("J96").Value = 1999.6666
$ws.Range("K96").Value = 7875.75
$ws.Range("L96").Value = 1999.6666
$ws.Range("M96").Value = -6502.75
$ws.Range("N96").Value = -4745.6666
$ws.Range("H107").Value = 1425.9
$ws.Range("I107").Value = 790.7692
$ws.Range("J107").Value = 2605.4285
$ws.Range("K107").Value = 2372.3076
$ws.Range("L107").Value = 7816.2855
$ws.Range("M107").Value = -452.3076000000001
$ws.Range("N107").Value = -11656.2855
$ws.Range("H122").Value = 3203.4443
$ws.Range("I122").Value = 3648.8572
$ws.Range("K122").Value = 10946.5716
$ws.Range("M122").Value = -8496.571599999999
$ws.Range("H132").Value = 2830.8333
$ws.Range("I132").Value = 1997
$ws.Range("K132").Value = 5991
$ws.Range("M132").Value = -3461
$ws.Range("H136").Value = 882.1667
$ws.Range("I136").Value = 882.1667
$ws.Range("K136").Value = 2646.5001
$ws.Range("M136").Value = -96.5001000000002
